# Auto-generated edit script: updates market-price derived columns
# (currentAveragePrice* / LevePrice* / LeveProfit*) across all 8 sheets,
# matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H17" = 570.0270400000001
    "J17" = 570.0270400000001
    "L17" = 1710.08112
    "N17" = -2046.08112
    "H43" = 1499.3334
    "I43" = 1499.3334
    "K43" = 1499.3334
    "M43" = -1430.3334
    "H57" = 59969
    "J57" = 59969
    "L57" = 179907
    "N57" = -180905
    "H92" = 3503.8
    "I92" = 1479.5714
    "J92" = 5275
    "K92" = 1479.5714
    "L92" = 5275
    "M92" = -231.5714
    "N92" = -7771
    "H98" = 1659.375
    "I98" = 1182.1428
    "K98" = 1182.1428
    "M98" = 315.8571999999999
    "H99" = 385.0625
    "I99" = 307.66666
    "K99" = 922.9999799999999
    "M99" = 575.0000200000001
    "H122" = 1659.375
    "I122" = 1182.1428
    "K122" = 3546.4284
    "M122" = -1096.4284
    "H132" = 1672333.4
    "I132" = 6666.6665
    "K132" = 19999.9995
    "M132" = -17469.9995
    "H137" = 1196.5128
    "I137" = 768.9286
    "J137" = 1435.96
    "K137" = 2306.7858
    "L137" = 4307.88
    "M137" = 243.2142000000003
    "N137" = -9407.880000000001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H32" = 4433.549
    "I32" = 1691.3556
    "K32" = 1691.3556
    "M32" = -1404.3556
    "H45" = 3531.375
    "I45" = 2621.75
    "K45" = 2621.75
    "M45" = -2244.75
    "H61" = 3558.25
    "I61" = 2621.25
    "K61" = 2621.25
    "M61" = -2409.25
    "H132" = 2312.3
    "I132" = 2347
    "K132" = 7041
    "M132" = -4511
    "H136" = 3558.25
    "I136" = 2621.25
    "K136" = 7863.75
    "M136" = -5313.75
    "H137" = 0
    "J137" = 0
    "L137" = 0
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
foreach ($cell in @("N137")) {
    $ws.Range($cell).ClearContents()
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H25" = 4999.5
    "J25" = 0
    "L25" = 0
    "H26" = 17499.5
    "I26" = 17499.5
    "K26" = 17499.5
    "M26" = -17207.5
    "H86" = 2579.889
    "I86" = 1370.3334
    "K86" = 1370.3334
    "M86" = -247.3334
    "H89" = 2579.889
    "I89" = 1370.3334
    "K89" = 6851.666999999999
    "M89" = -1235.666999999999
    "H99" = 4056.2856
    "I99" = 4119.4
    "K99" = 4119.4
    "M99" = -2621.4
    "H105" = 2758.7273
    "I105" = 2677.25
    "J105" = 2805.2856
    "K105" = 2677.25
    "L105" = 2805.2856
    "M105" = -930.25
    "N105" = -6299.2856
    "H134" = 2701
    "I134" = 2701
    "J134" = 0
    "K134" = 8103
    "L134" = 0
    "M134" = -5568
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
foreach ($cell in @("N25", "N134")) {
    $ws.Range($cell).ClearContents()
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 2288.88
    "I31" = 1765.2858
    "J31" = 2955.2727
    "K31" = 1765.2858
    "L31" = 2955.2727
    "M31" = -1470.2858
    "N31" = -3545.2727
    "H34" = 2288.88
    "I34" = 1765.2858
    "J34" = 2955.2727
    "K34" = 1765.2858
    "L34" = 2955.2727
    "M34" = -1563.2858
    "N34" = -3359.2727
    "H58" = 2160.9167
    "J58" = 2094.5
    "L58" = 2094.5
    "N58" = -2500.5
    "H99" = 37489
    "I99" = 8047.625
    "J99" = 115999.336
    "K99" = 8047.625
    "L99" = 115999.336
    "M99" = -6549.625
    "N99" = -118995.336
    "H122" = 2637.3333
    "I122" = 2956
    "J122" = 2000
    "K122" = 8868
    "L122" = 6000
    "M122" = -6418
    "N122" = -10900
    "H126" = 37489
    "I126" = 8047.625
    "J126" = 115999.336
    "K126" = 24142.875
    "L126" = 347998.008
    "M126" = -21672.875
    "N126" = -352938.008
    "H132" = 0
    "I132" = 0
    "K132" = 0
    "H134" = 3585.1
    "I134" = 3650.111
    "J134" = 3000
    "K134" = 10950.333
    "L134" = 9000
    "M134" = -8415.332999999999
    "N134" = -14070
    "H136" = 2160.9167
    "J136" = 2094.5
    "L136" = 6283.5
    "N136" = -11383.5
    "H137" = 110000
    "J137" = 110000
    "L137" = 110000
    "N137" = -120200
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
foreach ($cell in @("M132")) {
    $ws.Range($cell).ClearContents()
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H122" = 360.33334
    "J122" = 392.4
    "L122" = 3531.6
    "N122" = -8431.6
    "H137" = 3799.8
    "I137" = 2000
    "J137" = 4249.75
    "K137" = 6000
    "L137" = 12749.25
    "M137" = -900
    "N137" = -22949.25
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H102" = 4140
    "J102" = 4140
    "L102" = 4140
    "N102" = -7384
    "H126" = 2752.5
    "I126" = 2255.5
    "J126" = 3249.5
    "K126" = 6766.5
    "L126" = 9748.5
    "M126" = -4296.5
    "N126" = -14688.5
    "H132" = 6165.3
    "I132" = 6081.625
    "J132" = 6500
    "K132" = 18244.875
    "L132" = 19500
    "M132" = -15714.875
    "N132" = -24560
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H16" = 318.13333
    "I16" = 318.13333
    "K16" = 318.13333
    "M16" = -148.13333
    "H122" = 3360.5417
    "I122" = 3336.647
    "K122" = 10009.941
    "M122" = -7559.940999999999
    "H132" = 3333.6
    "I132" = 2826.353
    "J132" = 4411.5
    "K132" = 8479.059000000001
    "L132" = 13234.5
    "M132" = -5949.059000000001
    "N132" = -18294.5
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H122" = 5120.4185
    "I122" = 5691.6787
    "J122" = 4054.0667
    "K122" = 17075.0361
    "L122" = 12162.2001
    "M122" = -14625.0361
    "N122" = -17062.2001
}
foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Host "Applied $([string]205) cell updates and $([string]4) clears across 8 sheets."